$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend the numbered header by two more columns ---
# P1/Q1 need the same bold / bordered / centered format as the existing
# header cells (e.g. O1), so copy formats from O1 first, then set values.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Data rows 2-25 ---
# For each row: columns I, K, M, O swap their 1/2 values, and two new
# columns P, Q are appended with value 2 (plain format, like the other
# data cells in the row).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2  # column I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # column K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # column M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # column O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # column P: new
    $ws.Cells.Item($r, 17).Value = 2  # column Q: new
}

Write-Output "edit applied"
